$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.981.74"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "1.651.30"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "309.88"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "0.3904"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("D8").Value = "0.3818"
$ws.Range("E8").Value = "  -2.48%  "
$ws.Range("D9").Value = "52.05"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("E10").Value = "  -4.30%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "0.08453"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "23.87"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").Value = "7.063"
$ws.Range("E14").Value = "  -3.54%  "
$ws.Range("D15").Value = "8.009"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("D16").Value = "0.00001310"
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("D17").Value = "1.656.49"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").Value = "94.49"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").Value = "0.07007"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("D20").Value = "19.71"
$ws.Range("E20").Value = "  -4.12%  "
$ws.Range("D21").Value = "6.984"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").Value = "13.80"
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").Value = "23.974.95"
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("D25").Value = "2.444"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").Value = "2.928"
$ws.Range("E26").Value = "  -2.68%  "
$ws.Range("D27").Value = "22.08"
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("D28").Value = "152.98"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("D29").Value = "5.418"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("D30").Value = "138.10"
$ws.Range("E30").Value = "  -3.12%  "
$ws.Range("D31").Value = "7.934"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").Value = "2.536"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").Value = "1.829.23"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").Value = "1.019"
$ws.Range("E34").Value = "  -4.43%  "
$ws.Range("D35").Value = "0.08066"
$ws.Range("E35").Value = "  -2.18%  "
$ws.Range("D36").Value = "6.729"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").Value = "0.02925"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("D38").Value = "0.2680"
$ws.Range("E38").Value = "  -3.30%  "
$ws.Range("D39").Value = "10.74"
$ws.Range("E39").Value = "  -3.34%  "
$ws.Range("D40").Value = "0.09113"
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("D41").Value = "0.7605"
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("D42").Value = "13.40"
$ws.Range("E42").Value = "  -3.26%  "
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("D44").Value = "16.34"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D45").Value = "0.6968"
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").Value = "2.458"
$ws.Range("E46").Value = "  -3.15%  "
$ws.Range("D47").Value = "4.099"
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").Value = "0.08335"
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("D50").Value = "134.88"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("D51").Value = "1.225"
$ws.Range("E51").Value = "  -3.34%  "
